# Roteiro.xlsx — "Add files via upload" edit
#
# Sheet1 "1. Trade de GAP em Cripto": fill in the "Definição" column (C) for
# the strategic items, rename a couple of items, and add two new strategic
# items ("Definir a corretora à usar" / "Definir roteiro de desenvolvimento
# (etapas)") before the single "Desenvolver" dev item.
#
# Sheet2 "2. Opção binária": untouched content-wise (only shared-string
# indices shift upstream because unused strings were removed — nothing to do
# here, COM writes re-derive the shared-string table automatically).
#
# Sheet3 "3. Rede social profissional": the Item/Definição work hadn't
# started yet — fill in the "Item" column (B) that was previously blank, and
# the trailing Setor/Status cells that were missing.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "1. Trade de GAP em Cripto"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("1. Trade de GAP em Cripto")

# B/item renames for the rows that already existed
$ws1.Range("B2").Value = "Definir o problema"
$ws1.Range("B3").Value = "Definir o propósito"
$ws1.Range("B4").Value = "Definir o público alvo (país/público)"
$ws1.Range("B5").Value = "Definir a estratégia"
$ws1.Range("B8").Value = "Desenvolver"

# New rows 6/7, and the new "Definição" column — entered in the same order
# the author originally typed them (matches the shared-string table order):
# B6, B7, C2, C3, C4, C6 (Binance), then C5 (Shortar) last.
$ws1.Range("B6").Value = "Definir a corretora à usar"
$ws1.Range("B7").Value = "Definir roteiro de desenvolvimento (etapas)"

$ws1.Range("C2").Value = "Volatilidade alta no mercado de criptomoedas que dá margem para ganhos entre a cotação do contrato futuro de alguma moeda e o seu spot"
$ws1.Range("C3").Value = "Se aproveitar dessa volatilidade para lucros rápidos e consistentes até que o gap se estabilize igual no mercado ""padrão"" (em que este gap é ínfimo)"
$ws1.Range("C4").Value = "Podemos iniciar com nosso próprio capital e deixar rodando por alguns meses. Caso continue lucrativo, pensar em captar capital de terceiros para alavancar (?)"
$ws1.Range("C6").Value = "Binance (?) (maior do mundo aparentemente)"
$ws1.Range("C5").Value = "Shortar no maior e longar no menor entre cripto spot e mesma cripto futura com prazo (tem que tomar cuidado com o funding contrário à nossa posição do futuro)"

$ws1.Activate()
$ws1.Range("B19").Select()

# ---------------------------------------------------------------------
# Sheet 2: "2. Opção binária" — content unchanged, just re-select.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("2. Opção binária")
$ws2.Activate()
$ws2.Range("B34").Select()

# ---------------------------------------------------------------------
# Sheet 3: "3. Rede social profissional"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("3. Rede social profissional")

# Every row's "Setor" becomes "Estratégico" (it used to cycle through
# Desenvolvimento/Financeiro/Estratégico/Marketing/blank/blank — the
# strategy phase items hadn't been entered yet) and the "Item" column gets
# populated for the first time.
$ws3.Range("A2").Value = "Estratégico"
$ws3.Range("B2").Value = "Definir o problema"
$ws3.Range("E2").Value = "Em andamento"

$ws3.Range("A3").Value = "Estratégico"
$ws3.Range("B3").Value = "Definir o propósito"
$ws3.Range("E3").Value = "À iniciar"

$ws3.Range("A4").Value = "Estratégico"
$ws3.Range("B4").Value = "Definir o público alvo (país/público)"
$ws3.Range("E4").Value = "À iniciar"

$ws3.Range("A5").Value = "Estratégico"
$ws3.Range("B5").Value = "Definir a estratégia"
$ws3.Range("E5").Value = "À iniciar"

$ws3.Range("A6").Value = "Estratégico"
$ws3.Range("B6").Value = "Definir as tecnologias à utilizar"
$ws3.Range("E6").Value = "À iniciar"

$ws3.Range("A7").Value = "Estratégico"
$ws3.Range("B7").Value = "Definir freelancer ou sócio desenvolvedor"
$ws3.Range("E7").Value = "À iniciar"

$ws3.Activate()
$ws3.Range("C14").Select()

# ---------------------------------------------------------------------
# Best-effort column width refresh (Excel auto bestFit reacts to the new,
# longer text dropped into columns B/C on sheet 1 and B on sheet 3).
# ---------------------------------------------------------------------
$ws1.Columns.Item(2).AutoFit()
$ws1.Columns.Item(3).AutoFit()
$ws3.Columns.Item(2).AutoFit()

# Re-activate sheet 1 as the visible tab (tabSelected="1" in sheet1.xml).
$ws1.Activate()

$wb.Save()
